# Update the date line and regenerate the practice-sheet answers to match
# the output generated at c8c62b6.
#
# Word's Find/Replace matches exact text (MatchCase + MatchWholeWord), and
# every "old" value below is unique in the document, so a straightforward
# sequence of Find.Execute replacements is sufficient. The only subtlety is
# ordering: two of the new answers happen to equal an *old* answer that
# itself still needs to be replaced later (e.g. "27÷6=4, 3" becomes
# "87÷2=43, 1", while the original "87÷2=43, 1" cell becomes "58÷2=29, 0").
# To avoid a later replacement accidentally re-matching text that was just
# written, every "terminal" replacement in such a chain is performed before
# the replacement that produces that text.

$d = $word.ActiveDocument

function ReplaceExact($find, $repl) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $repl, 2) | Out-Null
}

# Header date.
ReplaceExact "2025-07-27 Sunday" "2025-07-28 Monday"

# Row 1 (first data row of the table).
ReplaceExact "36÷6=6, 0"   "64÷8=8, 0"
ReplaceExact "94÷7=13, 3"  "44÷5=8, 4"
ReplaceExact "65÷3=21, 2"  "77÷5=15, 2"
ReplaceExact "58÷7=8, 2"   "54÷7=7, 5"

# Row 2.
ReplaceExact "50÷4=12, 2"  "55÷3=18, 1"
ReplaceExact "87÷2=43, 1"  "58÷2=29, 0"
ReplaceExact "92÷2=46, 0"  "22÷9=2, 4"
ReplaceExact "27÷6=4, 3"   "87÷2=43, 1"
ReplaceExact "44÷2=22, 0"  "41÷8=5, 1"

# Row 3.
ReplaceExact "29÷5=5, 4"   "11÷5=2, 1"
ReplaceExact "86÷3=28, 2"  "39÷3=13, 0"
ReplaceExact "33÷3=11, 0"  "71÷2=35, 1"
ReplaceExact "76÷3=25, 1"  "67÷6=11, 1"
ReplaceExact "43÷5=8, 3"   "69÷3=23, 0"

# Row 4.
ReplaceExact "96÷2=48, 0"  "19÷4=4, 3"
ReplaceExact "45÷5=9, 0"   "10÷5=2, 0"
ReplaceExact "46÷4=11, 2"  "86÷3=28, 2"
ReplaceExact "89÷8=11, 1"  "31÷2=15, 1"

# Row 5.
ReplaceExact "86÷2=43, 0"  "16÷8=2, 0"
ReplaceExact "31÷6=5, 1"   "22÷8=2, 6"
ReplaceExact "89÷3=29, 2"  "11÷6=1, 5"
ReplaceExact "49÷8=6, 1"   "88÷8=11, 0"
ReplaceExact "17÷2=8, 1"   "90÷7=12, 6"
